$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a header row (row 1 was previously unused/empty - data started at row 2).
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "SourceLanguage"
$ws.Range("C1").Value = "TargetLanguage"
$ws.Range("D1").Value = "Status"

# Match the blank-but-present E column cell that every data row already has.
$ws.Range("E2").Copy($ws.Range("E1"))

# Remove the "I_Am_Old" resource row (originally row 5).
$ws.Rows.Item(5).Delete()
